$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (string) even when it looks numeric,
# without leaving a lasting number-format style on the cell.
function Set-TextValue($rangeAddr, $text) {
    $rng = $ws.Range($rangeAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Row 2
$ws.Range("B2").Value = "Ehhd"
$ws.Range("C2").Value = "м^3"
Set-TextValue "D2" "88"
Set-TextValue "F2" "6"
Set-TextValue "G2" "528"
$ws.Range("I2").Value = "24.04.2021"

# Row 3
$ws.Range("B3").Value = "fwe"
$ws.Range("C3").Value = "кг"
Set-TextValue "D3" "1"
Set-TextValue "F3" "400"
Set-TextValue "G3" "400"
$ws.Range("I3").Value = "26.04.2021"

# Remove row 4 entirely (also updates the used-range dimension)
$ws.Rows("4:4").Delete()
